# Quarterly income statement update:
#  - drop the oldest quarter (column D) so every quarter shifts one column left
#  - append the new quarter (column M): "فصل چهارم منتهی به 1401/12"
#  - the quarter that is now in column I ("فصل چهارم منتهی به 1400/12") was
#    re-published on 1402-02-27, which came with revised Sales / COGS figures
#    (read_price algorithm change) for that same quarter
#  - widen the new last column to match the "تاریخ انتشار" row's column style

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the oldest quarter - shifts E:M -> D:L
$ws.Columns.Item(4).Delete()

# 2) Clone column L's formatting (fills/borders/fonts/number formats) into the
#    freshly exposed column M so the new quarter matches its neighbours
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Column widths: the "wide" (31) column tracks the تاریخ انتشار column,
#    which is now M
$ws.Columns.Item(13).ColumnWidth = 31

# 4) New quarter header + publish date
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-27"

# 5) New quarter figures
$ws.Range("M11").Value = 1911
$ws.Range("M12").Value = -1410
$ws.Range("M13").Value = 501
$ws.Range("M14").Value = -77
$ws.Range("M15").Value = "-"
$ws.Range("M16").Value = 56
$ws.Range("M17").Value = 480
$ws.Range("M18").Value = -291
$ws.Range("M19").Value = 13
$ws.Range("M20").Value = 202
$ws.Range("M21").Value = -15
$ws.Range("M22").Value = 187
$ws.Range("M23").Value = "-"
$ws.Range("M24").Value = 187
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = 1705
$ws.Range("M27").Value = 0

# 6) The quarter in column I ("فصل چهارم منتهی به 1400/12") was re-published
#    on 1402-02-27 (was 1401-10-28 (6)) together with restated Sales/COGS
$ws.Range("I9").Value = "1402-02-27 (7)"
$ws.Range("I11").Value = 2225
$ws.Range("I12").Value = -1526
